$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header formatting (bold, border, centered/top-aligned) from an existing
# header cell (H1) onto the two new header cells so they match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new columns I (I0) and J (IF)
$data = @{
    2  = @(6, 8)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(10, 10)
    6  = @(1, 6)
    7  = @(1, 4)
    8  = @(8, 9)
    9  = @(6, 7)
    10 = @(3, 4)
    11 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
